$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 198 (「身体の一部として…」 post), shifting all rows below up by one.
$ws.Rows.Item(198).Delete()
